$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full replacement data set for rows 2-10 (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster, then the 16 numeric metric columns E:T).
# This expands the previous 2x2 (ECs/sCs) cross table into a 3x3
# (ECs/FAPs/sCs) cross table per Dr Hou's advice, adding the FAPs cluster.

$rows = @(
    @(2,  "ECs",  "Vwf", "Gp1ba", "ECs",  3, 1, 58.356725,              175.070175,           0.9850455027887332,  0.9850455027887334,  3, 1, 2.942633666666667, 8.827901000000001, 0.4532993947193931, 0.4532993947193931, 171.7224636614084,  1545.502172952675, 0.4465205301851929,  0.4465205301851931),
    @(3,  "ECs",  "Vwf", "Gp1ba", "FAPs", 3, 1, 58.356725,              175.070175,           0.9850455027887332,  0.9850455027887334,  2, 0.6666666666666666, 1.465278,        4.395834,           0.2257194424231681, 0.2257194424231681, 85.50882529454999,  769.57942765095,   0.2223439216509221,  0.2223439216509222),
    @(4,  "ECs",  "Vwf", "Gp1ba", "sCs",  3, 1, 58.356725,              175.070175,           0.9850455027887332,  0.9850455027887334,  3, 1, 2.083678,           6.251034000000001, 0.3209811628574388, 0.3209811628574388, 121.59662403455,    1094.36961631095,  0.316181050952618,   0.3161810509526181),
    @(5,  "FAPs", "Vwf", "Gp1ba", "ECs",  1, 0.3333333333333333, 0.3120586666666667, 0.936176, 0.005267464653066949, 0.00526746465306695, 3, 1, 2.942633666666667, 8.827901000000001, 0.4532993947193931, 0.4532993947193931, 0.9182743385084445, 8.264469046576,    0.002387738538941046, 0.002387738538941046),
    @(6,  "FAPs", "Vwf", "Gp1ba", "FAPs", 1, 0.3333333333333333, 0.3120586666666667, 0.936176, 0.005267464653066949, 0.00526746465306695, 2, 0.6666666666666666, 1.465278, 4.395834,         0.2257194424231681, 0.2257194424231681, 0.4572526989759999, 4.115274290784,    0.001188969184474018, 0.001188969184474019),
    @(7,  "FAPs", "Vwf", "Gp1ba", "sCs",  1, 0.3333333333333333, 0.3120586666666667, 0.936176, 0.005267464653066949, 0.00526746465306695, 3, 1, 2.083678, 6.251034000000001,               0.3209811628574388, 0.3209811628574388, 0.6502297784426667, 5.852068005984001, 0.001690756929651885, 0.001690756929651885),
    @(8,  "sCs",  "Vwf", "Gp1ba", "ECs",  3, 1, 0.5738856666666666, 1.721657, 0.00968703255819983, 0.009687032558199832, 3, 1, 2.942633666666667, 8.827901000000001, 0.4532993947193931, 0.4532993947193931, 1.688735283550778, 15.198617551957, 0.004391125995259037, 0.004391125995259038),
    @(9,  "sCs",  "Vwf", "Gp1ba", "FAPs", 3, 1, 0.5738856666666666, 1.721657, 0.00968703255819983, 0.009687032558199832, 2, 0.6666666666666666, 1.465278, 4.395834,                         0.2257194424231681, 0.2257194424231681, 0.8409020418819999, 7.568118376938,   0.002186551587771941, 0.002186551587771942),
    @(10, "sCs",  "Vwf", "Gp1ba", "sCs",  3, 1, 0.5738856666666666, 1.721657, 0.00968703255819983, 0.009687032558199832, 3, 1, 2.083678, 6.251034000000001,                                 0.3209811628574388, 0.3209811628574388, 1.195792938148667, 10.762136443338,  0.003109354975168851, 0.003109354975168852)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $r[$col]
    }
}
